$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "Refine" label textbox: nudge its rotation/position slightly.
$textBox = $s.Shapes.Item("TextBox 5")
$textBox.Rotation = 343.63433837890625
$textBox.Left = 282.2952880859375
$textBox.Top = 146.53433227539062

# Circular arrow: enlarge it, move it, and reshape it (arrowhead/thickness/sweep).
$arrow = $s.Shapes.Item("Circular Arrow 65")
$arrow.Left = 278.0166931152344
$arrow.Top = 10.259370803833008
$arrow.Width = 343.2237854003906
$arrow.Height = 302.03631591796875

$arrow.Adjustments.Item(1) = 5085
